$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

# Update header row values stay the same - only data rows change.
# Row 2 becomes: Equip-001, DESC-7, PT-9, PDT-9, 7, 3, 3, 3, 2
$ws.Range("A2").Value = "Equip-001"
$ws.Range("B2").Value = "DESC-7"
$ws.Range("C2").Value = "PT-9"
$ws.Range("D2").Value = "PDT-9"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2

# Row 3 becomes: Equip-003, DESC-5, PT-5, PDT-6, 9, 7, 5, 2, 3
$ws.Range("A3").Value = "Equip-003"
$ws.Range("B3").Value = "DESC-5"
$ws.Range("C3").Value = "PT-5"
$ws.Range("D3").Value = "PDT-6"
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 3

# Rows 4 and 5 are removed entirely.
$ws.Rows("4:5").Delete()
